$wb = $excel.ActiveWorkbook

# --- "List" sheet: account.civility -> account.civility.label -------------
$wsList = $wb.Worksheets.Item("List")
$wsList.Range("F2").Value = '${account.civility.label}'

# --- "List" sheet: account_addressId -> account_homeAddress ---------------
$homeAddressHeader = '${msg.getProperty(''account_homeAddress'')}'
$wsList.Range("K1").Value = $homeAddressHeader
$wsList.Range("K2").Value = '${printer.print(account.homeAddress)}'

# --- "Search" sheet: insert a new row for the home address field ----------
$wsSearch = $wb.Worksheets.Item("Search")
$wsSearch.Rows.Item(14).Insert()
$wsSearch.Range("A14").Value = $homeAddressHeader
$wsSearch.Range("B14").Value = '${homeAddress}'
